$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "63.115.40"
Set-TextValue "E2" "  +0.39%  "
Set-TextValue "D3" "2.561.08"
Set-TextValue "E3" "  +1.16%  "
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "581.87"
Set-TextValue "E5" "  +2.34%  "
Set-TextValue "D6" "147.64"
Set-TextValue "E6" "  -0.30%  "
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "E8" "  +0.73%  "
Set-TextValue "E9" "  +1.68%  "
Set-TextValue "E10" "  -0.86%  "
Set-TextValue "E11" "  +0.24%  "
Set-TextValue "E12" "  +0.12%  "
Set-TextValue "D13" "27.55"
Set-TextValue "E13" "  -0.81%  "
Set-TextValue "D14" "3.021.21"
Set-TextValue "E14" "  +1.00%  "
Set-TextValue "D15" "63.038.38"
Set-TextValue "E15" "  +0.24%  "
Set-TextValue "E16" "  +2.07%  "
Set-TextValue "D17" "2.553.07"
Set-TextValue "E17" "  +0.63%  "
Set-TextValue "D18" "11.37"
Set-TextValue "E18" "  -1.70%  "
Set-TextValue "D19" "341.88"
Set-TextValue "E19" "  +1.56%  "
Set-TextValue "D20" "4.37"
Set-TextValue "E20" "  +1.91%  "
Set-TextValue "D21" "6.81"
Set-TextValue "E21" "  +1.00%  "
Set-TextValue "E22" "  +0.09%  "
Set-TextValue "D23" "65.92"
Set-TextValue "E23" "  +0.43%  "
Set-TextValue "D24" "2.680.59"
Set-TextValue "E24" "  +0.89%  "
Set-TextValue "E25" "  +2.51%  "
Set-TextValue "E26" "  +0.91%  "
Set-TextValue "E27" "  +0.04%  "
Set-TextValue "B28" "InternetComputer(DFINITY)"
Set-TextValue "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D28" "8.42"
Set-TextValue "E28" "  +0.82%  "
Set-TextValue "B29" "SuiNetwork"
Set-TextValue "C29" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D29" "1.47"
Set-TextValue "E29" "  -3.05%  "
Set-TextValue "D30" "7.89"
Set-TextValue "E30" "  +9.06%  "
Set-TextValue "E31" "  +6.10%  "
Set-TextValue "D32" "0.0₃0822"
Set-TextValue "E32" "  +1.53%  "
Set-TextValue "D33" "176.45"
Set-TextValue "E33" "  -0.82%  "
Set-TextValue "D34" "1.60"
Set-TextValue "E34" "  +0.96%  "
Set-TextValue "D35" "430.72"
Set-TextValue "E35" "  +4.27%  "
Set-TextValue "D36" "0.404"
Set-TextValue "E36" "  +1.01%  "
Set-TextValue "D37" "19.21"
Set-TextValue "E37" "  +1.98%  "
Set-TextValue "D38" "4.47"
Set-TextValue "E38" "  +1.59%  "
Set-TextValue "E39" "  +0.02%  "
Set-TextValue "E40" "  +0.43%  "
Set-TextValue "E41" "  +0.09%  "
Set-TextValue "D42" "39.68"
Set-TextValue "E42" "  +1.17%  "
Set-TextValue "D43" "151.88"
Set-TextValue "E43" "  -0.19%  "
Set-TextValue "D44" "3.81"
Set-TextValue "E44" "  +1.76%  "
Set-TextValue "D45" "21.02"
Set-TextValue "E45" "  +1.82%  "
Set-TextValue "D46" "0.0551"
Set-TextValue "E46" "  +5.92%  "
Set-TextValue "E47" "  +0.14%  "
Set-TextValue "E48" "  +0.92%  "
Set-TextValue "E49" "  +2.00%  "
Set-TextValue "D50" "18.38"
Set-TextValue "E50" "  -0.17%  "
Set-TextValue "D51" "1.71"
Set-TextValue "E51" "  -3.67%  "
